$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the various "ID_*" primary-key column labels scattered across the
# little per-table schema legend to the single unified "Identyfikator"
# label, and rename "Osoba_zarzadzajaca" (where it denotes a column name,
# not the table/section title) to "Nazwa_os_zarzadzajaca".
$ws.Range("A2").Value  = "Identyfikator"
$ws.Range("J2").Value  = "Identyfikator"
$ws.Range("G3").Value  = "Identyfikator"
$ws.Range("C5").Value  = "Identyfikator"
$ws.Range("E5").Value  = "Identyfikator"
$ws.Range("G8").Value  = "Identyfikator"
$ws.Range("A12").Value = "Nazwa_os_zarzadzajaca"
$ws.Range("C13").Value = "Identyfikator"
$ws.Range("G13").Value = "Identyfikator"
$ws.Range("E19").Value = "Identyfikator"
$ws.Range("G19").Value = "Identyfikator"
$ws.Range("E20").Value = "Nazwa_os_zarzadzajaca"
$ws.Range("G24").Value = "Identyfikator"

# This label was a duplicate now superseded by the unified column above it;
# drop the whole cell (not just its text).
$ws.Range("C8").Clear()

# Column E was marked best-fit at the old (narrower) text width; widen it
# manually now that it holds longer strings and drop the auto-fit flag.
$ws.Columns.Item(5).ColumnWidth = 21.6

# Leave the cursor where the author last left it.
$ws.Range("E20").Select() | Out-Null
